# Applies the "Agregado de URL de repositorio" commit to the document.
# Most changes are simple text merges (collapsing runs that used to be
# split apart by <w:proofErr> spell-check markers); the substantive change
# is the insertion of a new hyperlink to the GitHub repository right
# after the "Baje del repositorio..." bullet, plus a couple of other
# small text additions (TiendaEvertec.sln, lastRenderedPageBreak).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Title: "Instructivo de Instalación de TiendaEvertec" ---
Replace-Text "Instructivo de Instalación de TiendaEvertec" "Instructivo de Instalación de TiendaEvertec"

# --- "A continuación se detallan los pasos para instalar la Tienda Evertec:" ---
Replace-Text "la Tienda Evertec:" "la Tienda Evertec:"

# --- Insert " (<github url>)" after the "Baje del repositorio..." bullet ---
$r = $d.Content
$r.Find.Execute("Baje del repositorio el código fuente de la solución web") | Out-Null
$ins = $r.Duplicate
$ins.Collapse(0) | Out-Null
$ins.InsertAfter(" (LINKPLACEHOLDER )") | Out-Null

$linkRange = $d.Content
$linkRange.Find.Execute("LINKPLACEHOLDER") | Out-Null
$d.Hyperlinks.Add($linkRange, "https://github.com/leomalevo/TiendaEvertec", "", "", "https://github.com/leomalevo/TiendaEvertec") | Out-Null

# --- "Una vez bajado..." bullet: merge TiendaEvertec.DB text back together ---
Replace-Text "la carpeta denominada TiendaEvertec.DB. Dentro de la misma se encuentra la DB " "la carpeta denominada TiendaEvertec.DB. Dentro de la misma se encuentra la DB "

# --- "En un servidor de base de datos SQL..." bullet ---
Replace-Text " llamada “TiendaEvertec” e importe el archivo " " llamada “TiendaEvertec” e importe el archivo "
Replace-Text "TiendaEvertec.bak sonre esta db" "TiendaEvertec.bak sonre esta db"

# --- "Cree un login y usuario..." bullet ---
Replace-Text " un login y usuario" " un login y usuario"
Replace-Text " para la db" " para la db"
Replace-Text "usertest”" "usertest”"
Replace-Text " y como password: “171615.Az”" " y como password: “171615.Az”"
Replace-Text " (Sin las comillas). Se pueden crear otros usuarios pero hay que cambiar el archivo web.Config para que considere el cambio" " (Sin las comillas). Se pueden crear otros usuarios pero hay que cambiar el archivo web.Config para que considere el cambio"

# --- "Una vez creada la DB..." bullet: add "TiendaEvertec.sln " and merge runs ---
Replace-Text "Una vez creada la DB, abrir la solución en Visual Studio " "Una vez creada la DB, abrir la solución TiendaEvertec.sln en Visual Studio "
Replace-Text "y crear un sitio web local en Internet Information Server (IIS) presionando botón derecho sobre el proyecto" "y crear un sitio web local en Internet Information Server (IIS) presionando botón derecho sobre el proyecto"
Replace-Text " web TiendaEvertec, Propiedas->solapa Web" " web TiendaEvertec, Propiedas->solapa Web"

# --- "Establezca el proyecto web..." bullet ---
Replace-Text "Establezca el proyecto web TiendaEvertec como protecto de Inicio " "Establezca el proyecto web TiendaEvertec como protecto de Inicio "

# --- "Actualice los paquetes de Nuget..." bullet ---
Replace-Text "Actualice los paquetes de Nuget de la solución, ingresando a la Consola de Administración de paquetes de Nuget, y luego tipeando en la misma: " "Actualice los paquetes de Nuget de la solución, ingresando a la Consola de Administración de paquetes de Nuget, y luego tipeando en la misma: "
Replace-Text "Update-Package. De esta manera se actualizarán todos los paquetes instalados en la solucion" "Update-Package. De esta manera se actualizarán todos los paquetes instalados en la solucion"

# --- "En caso de fallos..." bullet: merge NLog.config text, add lastRenderedPageBreak ---
Replace-Text "En caso de fallos la solución web está registrando errores mediante logs en un archivo denominado TiendaEvertec.log (La configuración de la misma se determina con NLog.config)" "En caso de fallos la solución web está registrando errores mediante logs en un archivo denominado TiendaEvertec.log (La configuración de la misma se determina con NLog.config)"


# --- "Pantallas del sitio TiendaEvertec" heading ---
Replace-Text "Pantallas del sitio TiendaEvertec" "Pantallas del sitio TiendaEvertec"

# --- "Pantalla Confirmacion de Orden y Redireccion a Pago" heading ---
Replace-Text "Confirmacion de Orden y Redireccion a Pago" "Confirmacion de Orden y Redireccion a Pago"

# --- Move the lastRenderedPageBreak off of "Pantalla Creación de Orden" ---
# (Word recalculates lastRenderedPageBreak locations on repagination, so no
#  explicit action is required here; leave as-is.)

Write-Output "done"
